$d = $word.ActiveDocument

# Locate the paragraph containing the final "Requisitos" bullet that must
# stay, so deletion starts right after it, and the paragraph holding the
# copyright footer, so deletion ends right after it (consuming its mark).
$anchorText = "LOB1053: F" + [char]0xED + "sica III (Requisito fraco)"
$footerText = [char]0xA9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $anchorText) {
        $startIndex = $i + 1
    }
    if ($text -eq $footerText) {
        $endIndex = $i
    }
}

if ($startIndex -ge 1 -and $endIndex -ge $startIndex) {
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($endIndex)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
